$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add Grant's GitHub username in the previously-empty B4 cell
$ws.Range("B4").Value = "grantmd26"
